# Updated cryptos list -- refresh Price (D) and Volume(1h) (E) columns, rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.527.09"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.98"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.75"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4531"
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3584"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.33"
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07100"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8909"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07741"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.28"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.815.99"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.259"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.293"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.64"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008513"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.578.10"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.10"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.948"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.936"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.90"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.76"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.022"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.19"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.823"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08706"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7400"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.715"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.417"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.110"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.071"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01935"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.916"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05083"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5114"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.764"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.019"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4695"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.988"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.35"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05977"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.77"
$ws.Range("E51").Value = "  -1.41%  "
